$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.407.83'
$ws.Range("E2").Value = '  +0.98%  '

$ws.Range("D3").Value = '3.441.02'
$ws.Range("E3").Value = '  +0.59%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '413.99'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.26%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.71'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.19%  '

$ws.Range("E7").Value = '  -1.05%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.725'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.11%  '

$ws.Range("E10").Value = '  +1.17%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.73'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.21%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.47'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +3.65%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000218'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +7.24%  '

$ws.Range("D14").Value = '3.980.38'
$ws.Range("E14").Value = '  +0.49%  '

$ws.Range("E15").Value = '  -0.26%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.50'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -3.86%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '12.98'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +6.11%  '

$ws.Range("D18").Value = '3.461.69'
$ws.Range("E18").Value = '  +1.13%  '

$ws.Range("E19").Value = '  +0.00%  '

$ws.Range("D20").Value = '62.439.72'
$ws.Range("E20").Value = '  +0.91%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '475.74'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +7.71%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '91.05'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.42%  '

$ws.Range("E23").Value = '  +3.17%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.43'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +3.71%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.60'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +22.42%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.31'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +2.13%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '33.29'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.14%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.80'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.02%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.62'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '11.97'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.17%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.65'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.24%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.167'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.37%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.112'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.22%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '40.77'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -4.43%  '

$ws.Range("E35").Value = '  +0.11%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '58.54'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +9.63%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0490'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.62%  '

$ws.Range("E38").Value = '  +0.09%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.02'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +3.01%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.325'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +4.54%  '

$ws.Range("E41").Value = '  -0.25%  '

$ws.Range("E42").Value = '  -0.89%  '

$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.68'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +6.70%  '

$ws.Range("B44").Value = 'Monero'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '145.15'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.94%  '

$ws.Range("E45").Value = '  +3.44%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.07'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +4.34%  '

$ws.Range("E47").Value = '  +12.64%  '

$ws.Range("B48").Value = 'PEPE'
$ws.Range("C48").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D48").Value = '0.0₃0562'
$ws.Range("E48").Value = '  +38.24%  '

$ws.Range("B49").Value = 'Celestia'
$ws.Range("C49").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '16.41'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.91%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.43'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.37%  '

$ws.Range("B51").Value = 'BitcoinSV'
$ws.Range("C51").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '110.67'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +8.00%  '
